# Fruta / hortaliza, semanal
# Insert a new weekly record into the "Brócoli" sheet.
# A new row is inserted before the current row 534, shifting all subsequent
# rows down by one (old row 534 -> new row 535, ..., old row 552 -> new row 553),
# and the newly created row 534 is populated with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data down starting at row 534.
$ws.Rows(534).Insert()

$newRow = 534

$ws.Cells.Item($newRow, 1).Value = 5
$ws.Cells.Item($newRow, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item($newRow, 3).Value = 'Maule'
$ws.Cells.Item($newRow, 4).Value = 45075
$ws.Cells.Item($newRow, 5).Value = 7
$ws.Cells.Item($newRow, 6).Value = 100112023
$ws.Cells.Item($newRow, 7).Value = 'Brócoli'
$ws.Cells.Item($newRow, 8).Value = 'Sin especificar'
$ws.Cells.Item($newRow, 9).Value = 'Primera'
$ws.Cells.Item($newRow, 10).Value = 5000
$ws.Cells.Item($newRow, 11).Value = 500
$ws.Cells.Item($newRow, 12).Value = 500
$ws.Cells.Item($newRow, 13).Value = 500
$ws.Cells.Item($newRow, 14).Value = '$/unidad'
$ws.Cells.Item($newRow, 15).Value = 'Región del Maule'
$ws.Cells.Item($newRow, 16).Value = 500
$ws.Cells.Item($newRow, 17).Value = 1
$ws.Cells.Item($newRow, 18).Value = 'Hortaliza'
